# Fruta / hortaliza, semanal
# Two new weekly price-report rows are inserted at the top of the
# "Limón" block (rows 452-453), pushing the existing rows 452:553 down
# to 454:555.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 452, shifting everything below down by 2.
$ws.Rows("452:453").Insert()

# --- Row 452: new record -------------------------------------------------
$ws.Range("A452").Value = 11
$ws.Range("B452").Value = "Vega Monumental Concepción"
$ws.Range("C452").Value = "Bíobío"
$ws.Range("D452").Value = 44785
$ws.Range("E452").Value = 8
$ws.Range("F452").Value = "Fruta"
$ws.Range("G452").Value = 100102
$ws.Range("H452").Value = "Cítricos"
$ws.Range("I452").Value = 100102003
$ws.Range("J452").Value = "Limón"
$ws.Range("K452").Value = "Sin especificar"
$ws.Range("L452").Value = "1a amarillo"
$ws.Range("M452").Value = 400
$ws.Range("N452").Value = 4500
$ws.Range("O452").Value = 5000
$ws.Range("P452").Value = 4688
$ws.Range("Q452").Value = "$/malla 16 kilos"
$ws.Range("R452").Value = "Región de O'Higgins"
$ws.Range("S452").Value = 293
$ws.Range("T452").Value = 16

# --- Row 453: new record -------------------------------------------------
$ws.Range("A453").Value = 11
$ws.Range("B453").Value = "Vega Monumental Concepción"
$ws.Range("C453").Value = "Bíobío"
$ws.Range("D453").Value = 44785
$ws.Range("E453").Value = 8
$ws.Range("F453").Value = "Fruta"
$ws.Range("G453").Value = 100102
$ws.Range("H453").Value = "Cítricos"
$ws.Range("I453").Value = 100102003
$ws.Range("J453").Value = "Limón"
$ws.Range("K453").Value = "Sin especificar"
$ws.Range("L453").Value = "2a amarillo"
$ws.Range("M453").Value = 300
$ws.Range("N453").Value = 3500
$ws.Range("O453").Value = 4000
$ws.Range("P453").Value = 3750
$ws.Range("Q453").Value = "$/malla 16 kilos"
$ws.Range("R453").Value = "Región de O'Higgins"
$ws.Range("S453").Value = 234
$ws.Range("T453").Value = 16
